$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200; existing rows 200-289 shift down to 201-290,
# carrying the formatting of the row above (as Excel's native Insert does).
$ws.Rows("200:200").Insert()

# Populate the newly inserted row 200 with the new weekly record.
$ws.Range("A200").Value = 1
$ws.Range("B200").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C200").Value = "Arica y Parinacota"
$ws.Range("D200").Value = 44726
$ws.Range("E200").Value = 15
$ws.Range("F200").Value = 100114013
$ws.Range("G200").Value = "Zanahoria"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 90
$ws.Range("K200").Value = 12000
$ws.Range("L200").Value = 13000
$ws.Range("M200").Value = 12500
$ws.Range("N200").Value = "$/saco 25 kilos"
$ws.Range("O200").Value = "Valle de Camiña"
$ws.Range("P200").Value = 500
$ws.Range("Q200").Value = 25
$ws.Range("R200").Value = "Hortaliza"
